$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where column H has a confidence value (i.e. not a -100 "no match" row)
# get a new "x" marker in column J.
$rowsWithX = @(2, 3, 4, 7, 8, 9, 10, 13, 14, 15, 16, 19, 20)
foreach ($r in $rowsWithX) {
    $ws.Cells.Item($r, 10).Value = "x"
}

# Selection / view changes to match the saved workbook state
[void]$ws.Range("J21").Select()
